$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Order.classCode row (row 12): Binding Value Set, Min, Base Min
$elements.Range("Z12").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActClass"
# leading apostrophe forces these numeric-looking values to stay text (matches
# the source data, which stores Min/Base Min as "0"/"1" strings, not numbers)
$elements.Range("F12").Value = "'0"
$elements.Range("AG12").Value = "'0"

# Order.moodCode row (row 13): Binding Value Set
$elements.Range("Z13").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActMoodIntent"

# Column Z (Binding Value Set) widened to fit the new, longer content
$elements.Columns.Item(26).ColumnWidth = 52.8
